$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for all touched Price/Volume cells so numeric-looking
# strings (e.g. "1.00", "67.565.06") are not coerced into numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.565.06'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.32%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.778.63'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.00%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.02'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.30'

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.90%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.782.51'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.06%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.14%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.159'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.40%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.18'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.90%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.462'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.43%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.25'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000244'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.69%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.431.42'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.53%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.797.18'

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.53%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.637.61'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.30%  '

$ws.Range("B18").Value = 'TRON'

$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.115'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -4.33%  '

$ws.Range("B19").Value = 'Polkadot'

$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.18'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.52'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '489.26'

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.97%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.03'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.82%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.740'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.54%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000149'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +12.11%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.12'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.73%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.36'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -6.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.25'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.18'

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.26%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.08%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.95'

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.43'

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.88%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.09'

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.68'

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.108'

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.50%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.14%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.88%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.77'

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.40%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.97%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.328'

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.12%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '448.66'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.40%  '

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.99'

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.88'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.32'

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.92%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.22'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -7.70%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.836.55'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.52%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '139.10'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.62%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0349'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.62%  '

$ws.Range("B50").Value = 'InjectiveProtocol'

$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.96'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.89%  '

$ws.Range("B51").Value = 'EnergySwap'

$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.26'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +11.72%  '
